# edit.ps1 - apply the OOXML diff via Word COM-interop
#
# Changes required:
#   1) Run "edited"                              -> "edited "            (trailing space moves in)
#   2) Run " template with three variables: "     -> "template with three variables: " (leading space moves out)
#   3) Bookmark var1var01 -> var_1_C2CA6CE1E93F82A0B90EA47268CBB363
#   4) Bookmark var2var02 -> var_2_609C03481F72362290823F373B8CE3C2
#   5) Bookmark var3var03 -> var_3_11D9761583B8480C6FFEA105475D931E
#
# Note: this simulated Word engine merges adjacent runs that share identical
# run formatting whenever their text is mutated via Range.Text. The source
# document (edited in LibreOffice) keeps the three runs separate, so each
# text edit below is performed while the target run is briefly given a
# distinguishing format (Bold) to prevent the automatic merge, then the
# format is reverted on exactly the new range once the text is in place.

$d = $word.ActiveDocument

# --- 1) "edited" -> "edited " -------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("edited", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Bold = 1
$r1.Text = "edited "
$r1.Bold = 0

# --- 2) " template with three variables: " -> "template with three variables: " ---
$r2 = $d.Content
$r2.Find.Execute(" template with three variables: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Bold = 1
$r2.Text = "template with three variables: "
$r2.Bold = 0

# --- 3-5) rename the three bookmarks using md5-style hashed names -------
function Rename-Bookmark($doc, $oldName, $newName) {
    $bm = $doc.Bookmarks.Item($oldName)
    $rng = $bm.Range
    $bm.Delete()
    $doc.Bookmarks.Add($newName, $rng) | Out-Null
}

Rename-Bookmark $d "var1var01" "var_1_C2CA6CE1E93F82A0B90EA47268CBB363"
Rename-Bookmark $d "var2var02" "var_2_609C03481F72362290823F373B8CE3C2"
Rename-Bookmark $d "var3var03" "var_3_11D9761583B8480C6FFEA105475D931E"

Write-Output $d.Content.Text
foreach ($b in $d.Bookmarks) {
    Write-Output $b.Name
}
